$wb = $excel.ActiveWorkbook

# Sheet 1: "Realizar o Acesso ao site"
$ws1 = $wb.Worksheets.Item("Realizar o Acesso ao site")
$ws1.Range("C14").Value = "A tela se adequa aquilo que foi proposto"

# Sheet 2: "Realizar o Calculo do Site"
$ws2 = $wb.Worksheets.Item("Realizar o Calculo do Site")
$ws2.Range("C16").Value = "A palavra ""Peso"" deverá aparecer como indicador"
$ws2.Range("D17").Value = "Cursor de texto deverá ficar ativo e a palavra ""Peso"" desaparecerá"
$ws2.Range("B18").Value = "Colocar o "" Peso "" igual a 85kg e clicar em calcular "

# Sheet 3: "Verificar campos e botão"
$ws3 = $wb.Worksheets.Item("Verificar campos e botão")
$ws3.Range("B15").Value = "Apagar o campo de inserção de números e não inserir nenhum valor, clicando em calcular"
$ws3.Range("B16").Value = "Ao colocar o "" Peso "" em valor abaixo de 0 e clicar em calcular"
$ws3.Range("B17").Value = "Ao colocar o "" Peso "" em valor  como "" 700 ""kg e clicar em calcular"
$ws3.Range("B18").Value = "Colocar o "" Peso "" igual a 100klg e clicar em calcular "
